# Update the "Team of Outs" player table (A2:C19) to reflect the latest
# roster/position/team info.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("De'Aaron Fox",    "PG",          "San Antonio Spurs"),
    @("Shaedon Sharpe",  "SG,SF",       "Portland Trail Blazers"),
    @("Mikal Bridges",   "SG,SF,PF",    "New York Knicks"),
    @("Isaiah Collier",  "PG",          "Utah Jazz"),
    @("Harrison Barnes", "SF,PF",       "San Antonio Spurs"),
    @("Miles Bridges",   "SF,PF",       "Charlotte Hornets"),
    @("Brandon Clarke",  "PF,C",        "Memphis Grizzlies"),
    @("DeMar DeRozan",   "SF,PF",       "Sacramento Kings"),
    @("Brook Lopez",     "C",           "Milwaukee Bucks"),
    @("Nick Richards",   "C",           "Phoenix Suns"),
    @("Nikola Vucevic",  "PF,C",        "Chicago Bulls"),
    @("Tyler Herro",     "PG,SG",       "Miami Heat"),
    @("Scottie Barnes",  "PG,SG,SF,PF", "Toronto Raptors"),
    @("Josh Giddey",     "PG,SG,SF",    "Chicago Bulls"),
    @("Evan Mobley",     "PF,C",        "Cleveland Cavaliers"),
    @("Luka Doncic",     "PG,SG",       "Los Angeles Lakers"),
    @("Bobby Portis",    "PF,C",        "Milwaukee Bucks"),
    @("Ja Morant",       "PG",          "Memphis Grizzlies")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
